# [EI-979] Update the survey.xlsx data dictionary header labels.
#
# The "If_Condition" data-dictionary sheet (Sheet1) carries two header
# cells, I1/J1, that name the columns used for conditional branching.
# Relabel them from the "_Question" wording to the "_Goto" wording:
#   I1: "Then_Question" -> "Then_Goto"
#   J1: "Else_Question" -> "Else_Goto"
# Also move the sheet's active selection to I1 (it was previously on J1),
# matching where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# Leave the cursor on the cell that was just edited (was J1, now I1).
$null = $ws.Range("I1").Select()
